$p = $ppt.ActivePresentation
$p.Slides.Item($p.Slides.Count).Delete()
